$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C (Förändrad) for rows 2 through 32 from 45608 to 45609
for ($row = 2; $row -le 32; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45608) {
        $cell.Value2 = 45609
    }
}
